$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.24779303665832941
$ws.Range("A2").Value = -0.0059999999889406297
$ws.Range("A3").Value = -0.0039999999900377503
$ws.Range("A4").Value = -0.0079999999817861323
$ws.Range("A5").Value = -0.0029999999893028928
$ws.Range("A6").Value = -0.0070344481689268434
$ws.Range("A7").Value = -0.0034663864201136541
$ws.Range("A8").Value = -0.0099999999738829004
$ws.Range("A9").Value = -0.001999999987994272
$ws.Range("A10").Value = -0.0019999999881470387
$ws.Range("A11").Value = -0.0029999999863958848
$ws.Range("A12").Value = -0.0034999999855651609
$ws.Range("A13").Value = -0.0034999999861238251
$ws.Range("A14").Value = -0.0079999999784927667
$ws.Range("A15").Value = -0.00099999999118516314
$ws.Range("A16").Value = 0.032464233503575635
$ws.Range("A17").Value = -0.0019999999897608589
$ws.Range("A18").Value = -0.0039999999861910496
$ws.Range("A19").Value = -0.050828603997919775
$ws.Range("A20").Value = -0.0039999999916968676
$ws.Range("A21").Value = -0.0039999999916116025
$ws.Range("A22").Value = -0.0039999999915458773
$ws.Range("A23").Value = -0.0049999999877732293
$ws.Range("A24").Value = -0.019999999960449877
$ws.Range("A25").Value = -0.019999999959969372
$ws.Range("A26").Value = -0.0024999999867016953
$ws.Range("A27").Value = -0.0024999999860115807
$ws.Range("A28").Value = -0.0019999999838598015
$ws.Range("A29").Value = -0.006999999972808979
$ws.Range("A30").Value = -0.059999999878677546
$ws.Range("A31").Value = -0.0069999999708656446
$ws.Range("A32").Value = -0.0099999999654691862
$ws.Range("A33").Value = -0.0039999999759281479
